{"js": "const replacements = [\n  [\"483\u00d73=1449\", \"616\u00d78=4928\"],\n  [\"213\u00d79=1917\", \"563\u00d78=4504\"],\n  [\"404\u00d79=3636\", \"203\u00d79=1827\"],\n  [\"232\u00d79=2088\", \"710\u00d79=6390\"],\n  [\"292\u00d74=1168\", \"614\u00d72=1228\"],\n  [\"145\u00d77=1015\", \"260\u00d72=520\"],\n  [\"252\u00d76=1512\", \"455\u00d72=910\"],\n  [\"147\u00d76=882\", \"685\u00d79=6165\"],\n  [\"726\u00d78=5808\", \"158\u00d79=1422\"],\n  [\"119\u00d79=1071\", \"808\u00d72=1616\"],\n  [\"533\u00d75=2665\", \"257\u00d72=514\"],\n  [\"503\u00d77=3521\", \"249\u00d73=747\"],\n  [\"933\u00d77=6531\", \"479\u00d79=4311\"],\n  [\"757\u00d77=5299\", \"732\u00d74=2928\"],\n  [\"901\u00d73=2703\", \"116\u00d75=580\"],\n  [\"448\u00d78=3584\", \"693\u00d73=2079\"],\n  [\"521\u00d75=2605\", \"541\u00d78=4328\"],\n  [\"978\u00d77=6846\", \"188\u00d73=564\"],\n  [\"616\u00d73=1848\", \"432\u00d78=3456\"],\n  [\"249\u00d78=1992\", \"400\u00d76=2400\"],\n  [\"202\u00d78=1616\", \"634\u00d78=5072\"],\n  [\"678\u00d73=2034\", \"687\u00d72=1374\"],\n  [\"118\u00d73=354\", \"424\u00d74=1696\"],\n  [\"282\u00d74=1128\", \"140\u00d79=1260\"],\n  [\"245\u00d76=1470\", \"222\u00d77=1554\"],\n];\n\nfor (const [oldText, newText] of replacements) {\n  const searchResults = context.document.body.search(oldText, { matchCase: true, matchWholeWord: false });\n  searchResults.load(\"items\");\n  await context.sync();\n\n  for (const range of searchResults.items) {\n    range.insertText(newText, Word.InsertLocation.replace);\n  }\n  await context.sync();\n}\n", "ps1": "$d = $word.ActiveDocument\n\n$replacements = @(\n    @(\"483\u00d73=1449\", \"616\u00d78=4928\"),\n    @(\"213\u00d79=1917\", \"563\u00d78=4504\"),\n    @(\"404\u00d79=3636\", \"203\u00d79=1827\"),\n    @(\"232\u00d79=2088\", \"710\u00d79=6390\"),\n    @(\"292\u00d74=1168\", \"614\u00d72=1228\"),\n    @(\"145\u00d77=1015\", \"260\u00d72=520\"),\n    @(\"252\u00d76=1512\", \"455\u00d72=910\"),\n    @(\"147\u00d76=882\", \"685\u00d79=6165\"),\n    @(\"726\u00d78=5808\", \"158\u00d79=1422\"),\n    @(\"119\u00d79=1071\", \"808\u00d72=1616\"),\n    @(\"533\u00d75=2665\", \"257\u00d72=514\"),\n    @(\"503\u00d77=3521\", \"249\u00d73=747\"),\n    @(\"933\u00d77=6531\", \"479\u00d79=4311\"),\n    @(\"757\u00d77=5299\", \"732\u00d74=2928\"),\n    @(\"901\u00d73=2703\", \"116\u00d75=580\"),\n    @(\"448\u00d78=3584\", \"693\u00d73=2079\"),\n    @(\"521\u00d75=2605\", \"541\u00d78=4328\"),\n    @(\"978\u00d77=6846\", \"188\u00d73=564\"),\n    @(\"616\u00d73=1848\", \"432\u00d78=3456\"),\n    @(\"249\u00d78=1992\", \"400\u00d76=2400\"),\n    @(\"202\u00d78=1616\", \"634\u00d78=5072\"),\n    @(\"678\u00d73=2034\", \"687\u00d72=1374\"),\n    @(\"118\u00d73=354\", \"424\u00d74=1696\"),\n    @(\"282\u00d74=1128\", \"140\u00d79=1260\"),\n    @(\"245\u00d76=1470\", \"222\u00d77=1554\")\n)\n\nforeach ($pair in $replacements) {\n    $oldText = $pair[0]\n    $newText = $pair[1]\n\n    $rng = $d.Content\n    $find = $rng.Find\n    $find.ClearFormatting()\n    $find.Text = $oldText\n    $find.Replacement.ClearFormatting()\n    $find.Replacement.Text = $newText\n    $find.Forward = $true\n    $find.Wrap = 1\n    $find.MatchCase = $true\n    $find.MatchWholeWord = $false\n    $find.Execute([ref]$oldText, [ref]$false, [ref]$false, [ref]$false, [ref]$false, [ref]$false, [ref]$true, [ref]2, [ref]$false, [ref]$newText, [ref]2)\n}\n"}
